# Update the Anxa1-Dysf LR-pairs sheet with new TPM-derived values.
#
# The sheet has one row per (Sending cluster, Target cluster) pair. Columns
# G/H (ligand avg/total expression) are constant per Sending cluster, and
# columns K/L/M/N (receptor cells/rate/avg/total expression) are constant
# per Target cluster. Columns I/J/O/P are "derived specificity" values,
# i.e. each cluster's G/H/M/N value divided by the sum of that value across
# all six clusters. Columns Q/R/S/T are simple products: Q=G*M, R=H*N,
# S=I*O, T=J*P. This script updates the underlying raw values (G,H,K,L,M,N)
# with the new TPM numbers and recomputes all of the derived columns so the
# whole table stays internally consistent.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$clusters = @("ECs", "FAPs", "Inflammatory-Mac", "MuSCs", "Neutrophils", "Resolving-Mac")

# New raw ligand expression values per Sending cluster.
$G = @{
    "ECs"               = 31.084226
    "FAPs"               = 210.2874346666667
    "Inflammatory-Mac"  = 131.5521063333333
    "MuSCs"             = 126.3423615
    "Neutrophils"       = 120.618675
    "Resolving-Mac"     = 142.588219
}
$H = @{
    "ECs"               = 62.168452
    "FAPs"               = 630.862304
    "Inflammatory-Mac"  = 394.6563190000001
    "MuSCs"             = 252.684723
    "Neutrophils"       = 361.856025
    "Resolving-Mac"     = 427.7646569999999
}

# New raw receptor-expressing-cells / detection-rate / expression values per
# Target cluster.
$K = @{
    "ECs"               = 2
    "FAPs"               = 3
    "Inflammatory-Mac"  = 3
    "MuSCs"             = 2
    "Neutrophils"       = 3
    "Resolving-Mac"     = 1
}
$L = @{
    "ECs"               = 1
    "FAPs"               = 1
    "Inflammatory-Mac"  = 1
    "MuSCs"             = 1
    "Neutrophils"       = 1
    "Resolving-Mac"     = 0.3333333333333333
}
$M = @{
    "ECs"               = 51.196146
    "FAPs"               = 0.9364056666666668
    "Inflammatory-Mac"  = 0.02709666666666667
    "MuSCs"             = 5.8798265
    "Neutrophils"       = 0.07041833333333333
    "Resolving-Mac"     = 0.03308066666666667
}
$N = @{
    "ECs"               = 102.392292
    "FAPs"               = 2.809217
    "Inflammatory-Mac"  = 0.08129
    "MuSCs"             = 11.759653
    "Neutrophils"       = 0.211255
    "Resolving-Mac"     = 0.099242
}

# Totals across all clusters, used for the derived-specificity columns.
$sumG = 0
$sumH = 0
$sumM = 0
$sumN = 0
foreach ($c in $clusters) {
    $sumG = $sumG + $G[$c]
    $sumH = $sumH + $H[$c]
    $sumM = $sumM + $M[$c]
    $sumN = $sumN + $N[$c]
}

$I = @{}
$J = @{}
$O = @{}
$P = @{}
foreach ($c in $clusters) {
    $I[$c] = $G[$c] / $sumG
    $J[$c] = $H[$c] / $sumH
    $O[$c] = $M[$c] / $sumM
    $P[$c] = $N[$c] / $sumN
}

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $sending = $ws.Cells.Item($r, 1).Value2
    $target  = $ws.Cells.Item($r, 4).Value2

    if (-not ($clusters -contains $sending)) { continue }
    if (-not ($clusters -contains $target)) { continue }

    $gVal = $G[$sending]
    $hVal = $H[$sending]
    $iVal = $I[$sending]
    $jVal = $J[$sending]

    $kVal = $K[$target]
    $lVal = $L[$target]
    $mVal = $M[$target]
    $nVal = $N[$target]
    $oVal = $O[$target]
    $pVal = $P[$target]

    $qVal = $gVal * $mVal
    $rVal = $hVal * $nVal
    $sVal = $iVal * $oVal
    $tVal = $jVal * $pVal

    $ws.Cells.Item($r, 7).Value2  = $gVal   # G
    $ws.Cells.Item($r, 8).Value2  = $hVal   # H
    $ws.Cells.Item($r, 9).Value2  = $iVal   # I
    $ws.Cells.Item($r, 10).Value2 = $jVal   # J
    $ws.Cells.Item($r, 11).Value2 = $kVal   # K
    $ws.Cells.Item($r, 12).Value2 = $lVal   # L
    $ws.Cells.Item($r, 13).Value2 = $mVal   # M
    $ws.Cells.Item($r, 14).Value2 = $nVal   # N
    $ws.Cells.Item($r, 15).Value2 = $oVal   # O
    $ws.Cells.Item($r, 16).Value2 = $pVal   # P
    $ws.Cells.Item($r, 17).Value2 = $qVal   # Q
    $ws.Cells.Item($r, 18).Value2 = $rVal   # R
    $ws.Cells.Item($r, 19).Value2 = $sVal   # S
    $ws.Cells.Item($r, 20).Value2 = $tVal   # T
}
